# Update the confusion-matrix test fixture (confMatrix) with refreshed
# values from a re-run of the logistic-regression experiment.
#
# Sheet layout: row 1 / col A are header labels (A1 is blank), B1:F1 hold
# the predicted-class headers (1..5), and A2:A6 hold the actual-class
# labels (1..5). B2:F6 is the 5x5 confusion matrix body.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 is an empty label cell. Explicitly clear it so the workbook doesn't
# keep a stale/blank shared-string pointer around after the re-save.
$ws.Range("A1").ClearContents()

# Refreshed confusion-matrix counts (only the cells that actually moved).
$ws.Range("C2").Value = 2    # was 8
$ws.Range("B3").Value = 8    # was 2
$ws.Range("D3").Value = 3    # was 0
$ws.Range("E3").Value = 2    # was 0
$ws.Range("C4").Value = 0    # was 3
$ws.Range("E4").Value = 16   # was 21
$ws.Range("C5").Value = 0    # was 2
$ws.Range("D5").Value = 21   # was 16
$ws.Range("F5").Value = 11   # was 9
$ws.Range("E6").Value = 9    # was 11
